# 9th Stab - Cosmetic Changes
#
# The sheet tracks one "UN" marker column per report date. Two newer dates
# (Jun_15, Jun_17) need to be inserted ahead of the existing Jun_13 / Jun_10
# columns, so the layout becomes:
#   B = Jun_17, C = Jun_15, D = Jun_13 (was B), E = Jun_10 (was C)
# with every data row carrying the same "UN" marker across the two new
# columns as it already does for the existing ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of data rows (header row 1 + data rows 2..lastRow) before we
# touch anything.
$lastRow = $ws.UsedRange.Rows.Count

# Insert two new columns at B, shifting the old B (Jun_13) -> D and
# old C (Jun_10) -> E, along with all of their data.
$ws.Columns("B:C").Insert()

# Headers for the newly inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Mirror the existing "UN" marker down the two new columns for every
# data row (rows 2-lastRow).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Cosmetic: keep the new/shifted columns the same width (8 characters)
# as the original Jun_10 column.
$ws.Columns("C").ColumnWidth = 7.1667
$ws.Columns("D").ColumnWidth = 7.1667
$ws.Columns("E").ColumnWidth = 7.1667
